# Insert a new weekly price record for Brócoli at "Vega Monumental Concepción"
# right before the existing row 253 (date 2022-03-03), shifting the rest of
# the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(253).Insert()

$ws.Cells.Item(253, 1).Value = 11
$ws.Cells.Item(253, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(253, 3).Value = "Bíobío"
$ws.Cells.Item(253, 4).Value = 44722
$ws.Cells.Item(253, 5).Value = 8
$ws.Cells.Item(253, 6).Value = 100112023
$ws.Cells.Item(253, 7).Value = "Brócoli"
$ws.Cells.Item(253, 8).Value = "Sin especificar"
$ws.Cells.Item(253, 9).Value = "Primera"
$ws.Cells.Item(253, 10).Value = 2800
$ws.Cells.Item(253, 11).Value = 700
$ws.Cells.Item(253, 12).Value = 750
$ws.Cells.Item(253, 13).Value = 727
$ws.Cells.Item(253, 14).Value = "$/unidad"
$ws.Cells.Item(253, 15).Value = "Región Metropolitana"
$ws.Cells.Item(253, 16).Value = 727
$ws.Cells.Item(253, 17).Value = 1
$ws.Cells.Item(253, 18).Value = "Hortaliza"
